$d = $word.ActiveDocument

# 1. "Link to github: " -> "Link to GitHub: "
$d.Content.Find.Execute("Link to github: ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Link to GitHub: ", 2)

# 2. "...best fit line to more precisely forecast the output, which..."
#    -> "...best fit line to forecast the output more precisely, which..."
$d.Content.Find.Execute("to more precisely forecast the output, which", $true, $false, $false, $false, $false,
                         $true, 1, $false, "to forecast the output more precisely, which", 2)

# 3. "...has z-score of more than 2 and less than -2. "
#    -> "...has z-score of more than two and less than -2. "
$d.Content.Find.Execute("has z-score of more than 2 and less than -2.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "has z-score of more than two and less than -2.", 2)

# 4. "...There are 300153 entries in the data (i.e. 300153 row)."
#    -> "...There are 300153 entries in the data (i.e., 300153 row)."
$d.Content.Find.Execute("in the data (i.e. 300153 row)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "in the data (i.e., 300153 row)", 2)

# 5. "...to avoid biased results (i.e. getting good results by chance). "
#    -> "...to avoid biased results (i.e., getting good results by chance). "
$d.Content.Find.Execute("biased results (i.e. getting good results by chance)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "biased results (i.e., getting good results by chance)", 2)

# 6. "...standard scaling and minimum/maximum scaler. "
#    -> "...standard scaling, and minimum/maximum scaler. "
$d.Content.Find.Execute("standard scaling and minimum/maximum scaler", $true, $false, $false, $false, $false,
                         $true, 1, $false, "standard scaling, and minimum/maximum scaler", 2)

# 7. "...as a method of evaluation, however it wasn't as effective..."
#    -> "...as a method of evaluation; however, it wasn't as effective..."
$d.Content.Find.Execute("method of evaluation, however it wasn" + [char]0x2019 + "t as effective",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "method of evaluation; however, it wasn" + [char]0x2019 + "t as effective", 2)

# 7b. "...RMSE was. Also in hindsight the data was only..."
#     -> "...RMSE was. Also, in hindsight the data was only..."
$d.Content.Find.Execute("RMSE was. Also in hindsight the data was only", $true, $false, $false, $false, $false,
                         $true, 1, $false, "RMSE was. Also, in hindsight the data was only", 2)

# 7c. "...flights booked in a year so we haven't got data..."
#     -> "...flights booked in a year, so we haven't got data..."
$d.Content.Find.Execute("flights booked in a year so we haven" + [char]0x2019 + "t got data",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "flights booked in a year, so we haven" + [char]0x2019 + "t got data", 2)

# 8. "...outliers which could of caused incorrect learning"
#    -> "...outliers which could have caused incorrect learning"
$d.Content.Find.Execute("outliers which could of caused incorrect learning", $true, $false, $false, $false, $false,
                         $true, 1, $false, "outliers which could have caused incorrect learning", 2)

# 9. "...increasing number of iterations and changing the activation function)..."
#    -> "...increasing number of iterations, and changing the activation function)..."
$d.Content.Find.Execute("increasing number of iterations and changing the activation function",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "increasing number of iterations, and changing the activation function", 2)
